$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.556.96'
$ws.Range("E2").Value = '  +5.20%  '
$ws.Range("D3").Value = '3.466.47'
$ws.Range("E3").Value = '  +5.68%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''581.16'
$ws.Range("E5").Value = '  +5.99%  '
$ws.Range("D6").Value = '''157.85'
$ws.Range("E6").Value = '  +5.00%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").Value = '3.469.42'
$ws.Range("E8").Value = '  +5.52%  '
$ws.Range("D9").Value = '''0.547'
$ws.Range("E9").Value = '  +4.40%  '
$ws.Range("D10").Value = '''7.59'
$ws.Range("E10").Value = '  +1.65%  '
$ws.Range("E11").Value = '  +7.18%  '
$ws.Range("D12").Value = '''0.443'
$ws.Range("E12").Value = '  +1.86%  '
$ws.Range("D13").Value = '4.067.74'
$ws.Range("E13").Value = '  +5.84%  '
$ws.Range("E14").Value = '  -0.92%  '
$ws.Range("D15").Value = '''0.0000196'
$ws.Range("E15").Value = '  +9.59%  '
$ws.Range("D16").Value = '''27.65'
$ws.Range("E16").Value = '  +3.79%  '
$ws.Range("D17").Value = '64.580.77'
$ws.Range("E17").Value = '  +5.40%  '
$ws.Range("D18").Value = '3.464.16'
$ws.Range("E18").Value = '  +5.46%  '
$ws.Range("D19").Value = '''6.45'
$ws.Range("E19").Value = '  +1.14%  '
$ws.Range("D20").Value = '''14.37'
$ws.Range("E20").Value = '  +5.46%  '
$ws.Range("D21").Value = '''396.57'
$ws.Range("E21").Value = '  +4.49%  '
$ws.Range("D22").Value = '''8.54'
$ws.Range("E22").Value = '  +0.83%  '
$ws.Range("D23").Value = '''0.545'
$ws.Range("E23").Value = '  +2.23%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '''72.54'
$ws.Range("E24").Value = '  +3.00%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '''0.997'
$ws.Range("E25").Value = '  -0.33%  '
$ws.Range("D26").Value = '''0.0000119'
$ws.Range("E26").Value = '  +20.73%  '
$ws.Range("D27").Value = '''9.67'
$ws.Range("E27").Value = '  +10.52%  '
$ws.Range("D28").Value = '''0.180'
$ws.Range("E28").Value = '  +3.99%  '
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '''2.06'
$ws.Range("E30").Value = '  +5.18%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").Value = '''5.87'
$ws.Range("E31").Value = '  +7.49%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").Value = '''1.38'
$ws.Range("E32").Value = '  +8.49%  '
$ws.Range("B33").Value = 'RenderToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D33").Value = '''6.61'
$ws.Range("E33").Value = '  +5.53%  '
$ws.Range("D34").Value = '''23.72'
$ws.Range("E34").Value = '  +4.24%  '
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("D36").Value = '''6.93'
$ws.Range("E36").Value = '  +3.08%  '
$ws.Range("D37").Value = '''1.50'
$ws.Range("E37").Value = '  +3.15%  '
$ws.Range("D38").Value = '''159.86'
$ws.Range("E38").Value = '  +0.58%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '''0.0778'
$ws.Range("E39").Value = '  +6.93%  '
$ws.Range("D40").Value = '''1.89'
$ws.Range("E40").Value = '  +5.11%  '
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = '''27.81'
$ws.Range("E41").Value = '  +4.33%  '
$ws.Range("D42").Value = '2.874.84'
$ws.Range("E42").Value = '  +2.16%  '
$ws.Range("D43").Value = '''0.0324'
$ws.Range("E43").Value = '  +1.72%  '
$ws.Range("D44").Value = '''0.778'
$ws.Range("E44").Value = '  +5.45%  '
$ws.Range("D45").Value = '''4.42'
$ws.Range("E45").Value = '  +2.30%  '
$ws.Range("D46").Value = '''41.37'
$ws.Range("E46").Value = '  +2.54%  '
$ws.Range("D47").Value = '''1.10'
$ws.Range("E47").Value = '  +6.81%  '
$ws.Range("D48").Value = '''22.94'
$ws.Range("E48").Value = '  +5.11%  '
$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").Value = '''2.16'
$ws.Range("E49").Value = '  +20.46%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '''6.48'
$ws.Range("E50").Value = '  +3.27%  '
$ws.Range("B51").Value = 'SuiNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D51").Value = '''0.840'
$ws.Range("E51").Value = '  +3.68%  '
